# The deck originally carries the "Integral" (Red Violet) theme colours on
# its slide master (ppt/theme/theme1.xml). The commit swaps the deck over to
# the default "Office Theme" colour palette (the one that, before the edit,
# only lived in the unused ppt/theme/theme2.xml part).
#
# PowerPoint exposes the theme's 12 colour slots (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink) through Theme.ThemeColorScheme on the slide
# master - same order as the OOXML <a:clrScheme> children. Re-pointing every
# slot to the "Office" palette reproduces the colour change the commit made.

function RGBInt($r, $g, $b) {
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$cs = $p.SlideMaster.Theme.ThemeColorScheme

$cs.Item(1).RGB  = RGBInt 0x00 0x00 0x00   # dk1
$cs.Item(2).RGB  = RGBInt 0xFF 0xFF 0xFF   # lt1
$cs.Item(3).RGB  = RGBInt 0x44 0x54 0x6A   # dk2
$cs.Item(4).RGB  = RGBInt 0xE7 0xE6 0xE6   # lt2
$cs.Item(5).RGB  = RGBInt 0x5B 0x9B 0xD5   # accent1
$cs.Item(6).RGB  = RGBInt 0xED 0x7D 0x31   # accent2
$cs.Item(7).RGB  = RGBInt 0xA5 0xA5 0xA5   # accent3
$cs.Item(8).RGB  = RGBInt 0xFF 0xC0 0x00   # accent4
$cs.Item(9).RGB  = RGBInt 0x44 0x72 0xC4   # accent5
$cs.Item(10).RGB = RGBInt 0x70 0xAD 0x47   # accent6
$cs.Item(11).RGB = RGBInt 0x05 0x63 0xC1   # hlink
$cs.Item(12).RGB = RGBInt 0x95 0x4F 0x72   # folHlink
